# Elimina las dos primeras filas de datos (sismos del 2025-12-06),
# dejando el resto de los registros que se recorren dos filas hacia arriba.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2:F3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
